{"js": "// The document currently splits several headline paragraphs (the title,\n// the author line, and the abstract) into many single-word/single-space\n// runs. Collapse each of those paragraphs down to one run holding the\n// full paragraph text, leaving every other paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"Questions: Trigonometric identities (degrees)\",\n  \"Author\": \"Dzhemma Ruseva\",\n  \"Abstract\": \"A selection of questions on trigonometric identities, where angles are measured in degrees.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const newText = replacements[para.style];\n  if (newText !== undefined) {\n    para.getRange().insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The title, author and abstract paragraphs are each split across many\n# single-word/single-space runs. Collapse each paragraph down to a\n# single run holding the full paragraph text, leaving every other\n# paragraph (e.g. the \"Summary\" abstract title) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"Title\"    = \"Questions: Trigonometric identities (degrees)\"\n    \"Author\"   = \"Dzhemma Ruseva\"\n    \"Abstract\" = \"A selection of questions on trigonometric identities, where angles are measured in degrees.\"\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $paragraph = $d.Paragraphs.Item($i)\n    $styleName = $paragraph.Style.NameLocal\n    if ($replacements.ContainsKey($styleName)) {\n        $range = $paragraph.Range\n        $oldText = $range.Text.TrimEnd([char]13, [char]7)\n        $find = $range.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $replacements[$styleName], 2)\n    }\n}\n"}
